$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A13").Value = 23
$ws.Range("B13").Value = "2021-11-16T21:56:50-03:00"
$ws.Range("C13").Value = 5.6
$ws.Range("D13").Value = "A Culpa e das Estrelas"
$ws.Range("E13").Value = "Murilo jose"
